# R10 is being split out of the "R43 R29 R10 R1" 100k group into its own
# line (R10 is changing to 430k / 530k 0.6W 1% metal film, to increase the
# minimum release time). Insert a new row 23, duplicating the formatting of
# the existing row 23 ("R13 R4" / 56R), then fill in the new part's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 23 (copy + insert-copy) so the new row inherits the correct
# number formats / styles (left-aligned Value column, text Order-code column,
# currency Price/Total columns) instead of plain defaults.
$ws.Rows.Item(23).Copy()
$ws.Rows.Item(23).Insert()

# Overwrite the freshly inserted row 23 with the new R10 / 430k component.
$ws.Range("A23").Value = "R10"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = "430k"
$ws.Range("D23").Value = "Axial 6.8mm"
$ws.Range("E23").Value = "530k 0.6W 1% metal film"
$ws.Range("F23").Value = "Multicomp"
$ws.Range("G23").Value = "MCMF006FF4303A50"
$ws.Range("H23").Value = "Farnell"
$ws.Range("I23").Value = "2401823"
$ws.Range("J23").Value = 0.0207
$ws.Range("K23").Value = 0.0207

# The original "R43 R29 R10 R1" 100k row (now pushed down to row 38 by the
# insert above) loses R10: relabel its reference designators and drop its
# quantity from 4 to 3.
$ws.Range("A38").Value = "R43 R29 R1 "
$ws.Range("B38").Value = 3

# Leave the selection where the editor ended up after making the change.
[void]$ws.Range("B39").Select()
